$wb = $excel.ActiveWorkbook

# The "Overview" rollup sheet shares the same "Ready for handoff" string with
# the per-language sheets, so it needs to be updated too to keep them in sync.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column: "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # The handoff transform failed, so there is no more "Latest Handoff File" --
    # remove the hyperlink + cell content entirely.
    $ws.Range("C2").Hyperlinks.Delete()
    $ws.Range("C2").Clear()

    # Dates reset back to the zero-date sentinel, and the row is now ignored.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"
}
